$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Parameters")

# Update the eBirdRegion value (row 6, column B) to add Stanislaus County
$ws.Range("B6").Value = "US-CA-085,US-CA-099"

# Move the active selection to B6, matching the saved selection in the workbook
$ws.Range("B6").Select()
